# Apply the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All cell updates, keyed by A1 address -> new literal text.
$updates = [ordered]@{
    "D2" = "70.838.94"
    "E2" = "  +0.64%  "
    "D3" = "3.607.84"
    "E3" = "  -0.12%  "
    "E4" = "  +0.09%  "
    "D5" = "580.78"
    "E5" = "  -0.94%  "
    "D6" = "189.27"
    "E6" = "  -1.74%  "
    "D7" = "3.586.37"
    "E7" = "  -0.51%  "
    "D8" = "0.625"
    "E8" = "  -1.90%  "
    "E9" = "  +0.41%  "
    "D10" = "0.188"
    "E10" = "  +3.61%  "
    "D11" = "0.656"
    "E11" = "  -1.69%  "
    "D12" = "54.74"
    "E12" = "  -4.93%  "
    "D13" = "0.0000307"
    "E13" = "  +0.69%  "
    "D14" = "9.68"
    "E14" = "  -1.36%  "
    "D15" = "4.225.84"
    "E15" = "  +0.72%  "
    "D16" = "19.68"
    "E16" = "  -2.47%  "
    "D17" = "3.631.09"
    "E17" = "  +0.02%  "
    "D18" = "70.906.72"
    "E18" = "  +0.71%  "
    "D19" = "12.59"
    "E19" = "  -0.24%  "
    "E20" = "  -0.78%  "
    "D21" = "1.05"
    "E21" = "  +0.26%  "
    "D22" = "503.20"
    "E22" = "  +4.06%  "
    "D23" = "19.36"
    "E23" = "  +0.08%  "
    "D24" = "4.93"
    "E24" = "  -3.58%  "
    "D25" = "4.40"
    "E25" = "  -0.72%  "
    "D26" = "95.87"
    "E26" = "  +6.31%  "
    "D27" = "11.61"
    "E27" = "  +2.65%  "
    "D28" = "2.99"
    "E28" = "  -4.23%  "
    "D29" = "9.41"
    "E29" = "  -0.18%  "
    "D30" = "7.76"
    "E30" = "  -2.79%  "
    "D31" = "31.92"
    "E31" = "  -1.69%  "
    "D32" = "12.61"
    "E32" = "  +3.26%  "
    "D33" = "66.46"
    "E33" = "  -0.04%  "
    "D34" = "0.117"
    "E34" = "  -2.30%  "
    "D35" = "576.69"
    "E35" = "  -5.51%  "
    "D36" = "3.25"
    "E36" = "  +9.93%  "
    "D37" = "39.11"
    "E37" = "  -2.57%  "
    "D38" = "0.414"
    "E38" = "  +1.96%  "
    "B39" = "Dai"
    "C39" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D39" = "0.998"
    "E39" = "  +0.00%  "
    "B40" = "PEPE"
    "C40" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D40" = "0.0₃0800"
    "E40" = "  -3.74%  "
    "D41" = "3.32"
    "E41" = "  +5.02%  "
    "B42" = "Kaspa"
    "C42" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D42" = "0.138"
    "E42" = "  -6.46%  "
    "B43" = "Stacks"
    "C43" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D43" = "3.47"
    "E43" = "  -2.46%  "
    "D44" = "3.06"
    "E44" = "  -1.90%  "
    "D45" = "0.0457"
    "E45" = "  +0.80%  "
    "D46" = "3.55"
    "E46" = "  +4.87%  "
    "D47" = "3.220.10"
    "E47" = "  -2.74%  "
    "D48" = "9.58"
    "E48" = "  -0.98%  "
    "D49" = "0.136"
    "E49" = "  -1.36%  "
    "E50" = "  +28.54%  "
    "D51" = "1.00"
    "E51" = "  -0.11%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $isPriceNumberLike = ($addr.StartsWith("D")) -and ($value -match '^[+-]?\d+(\.\d+)?$')
    if ($isPriceNumberLike) {
        # This Price-column value parses as a plain number (e.g. "580.78").
        # Prefix with an apostrophe so Excel keeps storing it as literal text
        # (matching the source feed/export), rather than silently coercing it
        # into a numeric cell.
        $ws.Range($addr).Value = "'" + $value
    } else {
        $ws.Range($addr).Value = $value
    }
}
